$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 507
$ws1.Range("F6").Value = 487
$ws1.Range("F7").Value = 88
$ws1.Range("F10").Value = 6426
$ws1.Range("F11").Value = 218
$ws1.Range("F13").Value = 2587
$ws1.Range("F14").Value = 151
$ws1.Range("F15").Value = 255
$ws1.Range("F16").Value = 251
$ws1.Range("F17").Value = 496

# Sheet "全部类型" (fourth sheet) - same underlying rows, offset by 2
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 507
$ws4.Range("F8").Value = 487
$ws4.Range("F9").Value = 88
$ws4.Range("F13").Value = 6426
$ws4.Range("F15").Value = 218
$ws4.Range("F17").Value = 2587
$ws4.Range("F18").Value = 151
$ws4.Range("F19").Value = 255
$ws4.Range("F20").Value = 251
$ws4.Range("F21").Value = 496
